$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "45-15=30" "57+11=68"
Replace-Text "55-18=37" "20+36=56"
Replace-Text "65-24=41" "92-0=92"
Replace-Text "24+40=64" "68-24=44"
Replace-Text "58-28=30" "13+86=99"
Replace-Text "17+45=62" "6+34=40"
Replace-Text "31+55=86" "15+57=72"
Replace-Text "91-21=70" "11+13=24"
Replace-Text "19+79=98" "64+23=87"
Replace-Text "2+51=53" "98-76=22"
Replace-Text "99-95=4" "52-35=17"
Replace-Text "75-13=62" "26+8=34"
Replace-Text "10+4=14" "99-25=74"
Replace-Text "53-15=38" "7+43=50"
Replace-Text "16+4=20" "26-22=4"
Replace-Text "23+4=27" "8+31=39"
Replace-Text "60-44=16" "36-3=33"
Replace-Text "57+31=88" "15-3=12"
Replace-Text "76-2=74" "78+21=99"
Replace-Text "15+43=58" "0+33=33"
Replace-Text "54-49=5" "91-33=58"
Replace-Text "30+45=75" "66-59=7"
Replace-Text "15-15=0" "85-61=24"
Replace-Text "13+0=13" "22+4=26"
Replace-Text "83+2=85" "93-3=90"
Replace-Text "99-58=41" "81+6=87"
Replace-Text "95-72=23" "79-71=8"
Replace-Text "58-24=34" "31-9=22"
Replace-Text "55+30=85" "48+14=62"
Replace-Text "34-22=12" "45-8=37"
Replace-Text "10+34=44" "28+47=75"
Replace-Text "33-25=8" "6+82=88"
Replace-Text "87-77=10" "23-23=0"
Replace-Text "26+22=48" "76-46=30"
Replace-Text "35-26=9" "38+40=78"
Replace-Text "16+2=18" "9+53=62"
Replace-Text "86-26=60" "9+21=30"
Replace-Text "51-15=36" "80-16=64"
Replace-Text "24-23=1" "50+5=55"
Replace-Text "72-25=47" "31+58=89"
Replace-Text "12+71=83" "63+20=83"
Replace-Text "32+13=45" "15+16=31"
Replace-Text "39+14=53" "90-57=33"
Replace-Text "44+30=74" "50+20=70"
Replace-Text "40-7=33" "98-21=77"
Replace-Text "43+33=76" "65+16=81"
Replace-Text "36+8=44" "91-57=34"
Replace-Text "37+30=67" "53+23=76"
Replace-Text "35-24=11" "60-40=20"
Replace-Text "9+17=26" "32-10=22"
Replace-Text "90+9=99" "97-6=91"
Replace-Text "3+42=45" "31-13=18"
Replace-Text "71-69=2" "65-45=20"
Replace-Text "84-2=82" "67-51=16"
Replace-Text "0+28=28" "77-41=36"
Replace-Text "76-12=64" "73-29=44"
Replace-Text "7+55=62" "71+4=75"
Replace-Text "16+57=73" "69+12=81"
Replace-Text "8+75=83" "72-8=64"
Replace-Text "25-24=1" "55+39=94"
Replace-Text "71-29=42" "99-1=98"
Replace-Text "39+2=41" "27+20=47"
Replace-Text "76-70=6" "80-32=48"
Replace-Text "6-3=3" "78-31=47"
Replace-Text "76-0=76" "76-25=51"
Replace-Text "13+17=30" "33+51=84"
Replace-Text "21-19=2" "22-4=18"
Replace-Text "4+44=48" "72-46=26"
Replace-Text "20+71=91" "34+64=98"
Replace-Text "13+66=79" "27-24=3"
Replace-Text "47-41=6" "80-5=75"
Replace-Text "8+22=30" "4+32=36"
Replace-Text "91-31=60" "5+13=18"
Replace-Text "30+5=35" "38+10=48"
Replace-Text "44-5=39" "96-87=9"
Replace-Text "19+80=99" "58-3=55"
Replace-Text "71-3=68" "21+60=81"
Replace-Text "19-0=19" "32+23=55"
Replace-Text "99-35=64" "58-23=35"
Replace-Text "22-3=19" "58+6=64"
Replace-Text "57+13=70" "81+9=90"
Replace-Text "95-6=89" "8+87=95"
Replace-Text "27+32=59" "14+74=88"
Replace-Text "62-43=19" "15-10=5"
Replace-Text "95-37=58" "37-30=7"
Replace-Text "39+8=47" "35+39=74"
Replace-Text "31+28=59" "35+30=65"
Replace-Text "16+10=26" "90-72=18"
Replace-Text "71+1=72" "60+4=64"
Replace-Text "24+7=31" "19-5=14"
Replace-Text "21+14=35" "4+5=9"
Replace-Text "17+72=89" "68+1=69"
Replace-Text "18+1=19" "81+16=97"
Replace-Text "31-3=28" "92-29=63"
Replace-Text "21-17=4" "8+8=16"
Replace-Text "26+44=70" "9+36=45"
Replace-Text "14+40=54" "84-23=61"
Replace-Text "4+33=37" "71+24=95"
Replace-Text "13+73=86" "85-72=13"
Replace-Text "9+16=25" "53-49=4"
